# Scheduled runner refresh of Moogle market-price snapshots.
# Updates currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ
# (and the dependent LevePrice*/LeveProfit* columns) on each per-job sheet.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4638.45
$ws.Range("J17").Value = 4638.45
$ws.Range("L17").Value = 13915.35
$ws.Range("N17").Value = -14251.35
$ws.Range("H28").Value = 1024.8846
$ws.Range("I28").Value = 665.9545000000001
$ws.Range("J28").Value = 2999
$ws.Range("K28").Value = 665.9545000000001
$ws.Range("L28").Value = 2999
$ws.Range("M28").Value = -180.9545000000001
$ws.Range("N28").Value = -3969
$ws.Range("H43").Value = 2845
$ws.Range("H51").Value = 15580.117
$ws.Range("J51").Value = 13651.777
$ws.Range("L51").Value = 13651.777
$ws.Range("N51").Value = -14619.777
$ws.Range("H55").Value = 284.75
$ws.Range("I55").Value = 54
$ws.Range("K55").Value = 54
$ws.Range("M55").Value = 160
$ws.Range("H62").Value = 6043.769
$ws.Range("I62").Value = 4732
$ws.Range("K62").Value = 4732
$ws.Range("M62").Value = -4108
$ws.Range("H65").Value = 6043.769
$ws.Range("I65").Value = 4732
$ws.Range("K65").Value = 23660
$ws.Range("M65").Value = -20540
$ws.Range("H88").Value = 1482.5
$ws.Range("I88").Value = 1500.3334
$ws.Range("K88").Value = 1500.3334
$ws.Range("M88").Value = -1094.3334
$ws.Range("H91").Value = 1482.5
$ws.Range("I91").Value = 1500.3334
$ws.Range("K91").Value = 1500.3334
$ws.Range("M91").Value = -96.33339999999998
$ws.Range("H99").Value = 3670.6667
$ws.Range("J99").Value = 7181
$ws.Range("L99").Value = 21543
$ws.Range("N99").Value = -24539
$ws.Range("H116").Value = 17368.375
$ws.Range("I116").Value = 20117.285
$ws.Range("J116").Value = 15888.192
$ws.Range("K116").Value = 20117.285
$ws.Range("L116").Value = 15888.192
$ws.Range("M116").Value = -16675.285
$ws.Range("N116").Value = -22772.192

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1116.5
$ws.Range("I2").Value = 483.18182
$ws.Range("J2").Value = 3438.6667
$ws.Range("K2").Value = 483.18182
$ws.Range("L2").Value = 3438.6667
$ws.Range("M2").Value = -370.18182
$ws.Range("N2").Value = -3664.6667
$ws.Range("H32").Value = 10319.777
$ws.Range("I32").Value = 6860.3887
$ws.Range("K32").Value = 6860.3887
$ws.Range("M32").Value = -6573.3887
$ws.Range("H45").Value = 3848798.8
$ws.Range("I45").Value = 5265004.5
$ws.Range("K45").Value = 5265004.5
$ws.Range("M45").Value = -5264627.5
$ws.Range("H61").Value = 5004
$ws.Range("I61").Value = 4973.4043
$ws.Range("K61").Value = 4973.4043
$ws.Range("M61").Value = -4761.4043
$ws.Range("H116").Value = 1116.5
$ws.Range("I116").Value = 483.18182
$ws.Range("J116").Value = 3438.6667
$ws.Range("K116").Value = 483.18182
$ws.Range("L116").Value = 3438.6667
$ws.Range("M116").Value = 1810.81818
$ws.Range("N116").Value = -8026.6667
$ws.Range("H132").Value = 3200.8484
$ws.Range("I132").Value = 1911.7142
$ws.Range("J132").Value = 10420
$ws.Range("K132").Value = 5735.142599999999
$ws.Range("L132").Value = 31260
$ws.Range("M132").Value = -3205.142599999999
$ws.Range("N132").Value = -36320
$ws.Range("H136").Value = 5004
$ws.Range("I136").Value = 4973.4043
$ws.Range("K136").Value = 14920.2129
$ws.Range("M136").Value = -12370.2129

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1116.5
$ws.Range("I3").Value = 483.18182
$ws.Range("J3").Value = 3438.6667
$ws.Range("K3").Value = 483.18182
$ws.Range("L3").Value = 3438.6667
$ws.Range("M3").Value = -369.18182
$ws.Range("N3").Value = -3666.6667
$ws.Range("H22").Value = 671.2273
$ws.Range("I22").Value = 541.9375
$ws.Range("J22").Value = 1016
$ws.Range("K22").Value = 541.9375
$ws.Range("L22").Value = 1016
$ws.Range("M22").Value = -368.9375
$ws.Range("N22").Value = -1362
$ws.Range("H99").Value = 1239.9615
$ws.Range("I99").Value = 1170.2632
$ws.Range("K99").Value = 1170.2632
$ws.Range("M99").Value = 327.7367999999999
$ws.Range("H105").Value = 1084189.5
$ws.Range("I105").Value = 1431709.9
$ws.Range("J105").Value = 3015.111
$ws.Range("K105").Value = 1431709.9
$ws.Range("L105").Value = 3015.111
$ws.Range("M105").Value = -1429962.9
$ws.Range("N105").Value = -6509.111

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 23759.715
$ws.Range("J51").Value = 25332.834
$ws.Range("L51").Value = 25332.834
$ws.Range("N51").Value = -26804.834
$ws.Range("H61").Value = 23759.715
$ws.Range("J61").Value = 25332.834
$ws.Range("L61").Value = 25332.834
$ws.Range("N61").Value = -26028.834
$ws.Range("H93").Value = 7631.3335
$ws.Range("I93").Value = 7631.3335
$ws.Range("K93").Value = 7631.3335
$ws.Range("M93").Value = -5759.3335
$ws.Range("H122").Value = 3031.0386
$ws.Range("I122").Value = 2406.0625
$ws.Range("K122").Value = 7218.1875
$ws.Range("M122").Value = -4768.1875
$ws.Range("H132").Value = 4384.4727
$ws.Range("I132").Value = 2885.8696
$ws.Range("J132").Value = 12044
$ws.Range("K132").Value = 8657.6088
$ws.Range("L132").Value = 36132
$ws.Range("M132").Value = -6127.6088
$ws.Range("N132").Value = -41192
$ws.Range("H134").Value = 1918.0862
$ws.Range("I134").Value = 1499.6383
$ws.Range("K134").Value = 4498.9149
$ws.Range("M134").Value = -1963.9149

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 4259.8
$ws.Range("J55").Value = 6766.3335
$ws.Range("L55").Value = 20299.0005
$ws.Range("N55").Value = -20653.0005
$ws.Range("H140").Value = 1749.08
$ws.Range("I140").Value = 1281.75
$ws.Range("J140").Value = 1838.0952
$ws.Range("K140").Value = 3845.25
$ws.Range("L140").Value = 5514.2856
$ws.Range("M140").Value = 1334.75
$ws.Range("N140").Value = -15874.2856
$ws.Range("H141").Value = 5728.7896
$ws.Range("I141").Value = 4203.769
$ws.Range("J141").Value = 9033
$ws.Range("K141").Value = 12611.307
$ws.Range("L141").Value = 27099
$ws.Range("M141").Value = -7431.307000000001
$ws.Range("N141").Value = -37459

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 12305.9
$ws.Range("J57").Value = 41499
$ws.Range("L57").Value = 41499
$ws.Range("N57").Value = -43139
$ws.Range("H75").Value = 30000
$ws.Range("J75").Value = 30000
$ws.Range("L75").Value = 30000
$ws.Range("N75").Value = -31748
$ws.Range("H78").Value = 30000
$ws.Range("J78").Value = 30000
$ws.Range("L78").Value = 90000
$ws.Range("N78").Value = -98736
$ws.Range("H97").Value = 566.53845
$ws.Range("I97").Value = 566.4706
$ws.Range("J97").Value = 566.6667
$ws.Range("K97").Value = 566.4706
$ws.Range("L97").Value = 566.6667
$ws.Range("M97").Value = -70.47059999999999
$ws.Range("N97").Value = -1558.6667

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1159.4706
$ws.Range("J16").Value = 1413.7142
$ws.Range("L16").Value = 1413.7142
$ws.Range("N16").Value = -1753.7142
$ws.Range("H55").Value = 509.65384
$ws.Range("I55").Value = 290.93332
$ws.Range("J55").Value = 807.9091
$ws.Range("K55").Value = 290.93332
$ws.Range("L55").Value = 807.9091
$ws.Range("M55").Value = -117.93332
$ws.Range("N55").Value = -1153.9091
$ws.Range("H61").Value = 5677.2856
$ws.Range("I61").Value = 2110.7778
$ws.Range("K61").Value = 2110.7778
$ws.Range("M61").Value = -1908.7778
$ws.Range("H99").Value = 39321.25
$ws.Range("I99").Value = 29333.334
$ws.Range("K99").Value = 29333.334
$ws.Range("M99").Value = -26338.334
$ws.Range("H113").Value = 5677.2856
$ws.Range("I113").Value = 2110.7778
$ws.Range("K113").Value = 2110.7778
$ws.Range("M113").Value = 59.22220000000016
$ws.Range("H122").Value = 4060.64
$ws.Range("I122").Value = 2712.7856
$ws.Range("J122").Value = 5776.091
$ws.Range("K122").Value = 8138.3568
$ws.Range("L122").Value = 17328.273
$ws.Range("M122").Value = -5688.3568
$ws.Range("N122").Value = -22228.273
$ws.Range("H136").Value = 6942.029
$ws.Range("I136").Value = 4448.625
$ws.Range("K136").Value = 13345.875
$ws.Range("M136").Value = -10795.875

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H44").Value = 90225
$ws.Range("J44").Value = 90225
$ws.Range("L44").Value = 90225
$ws.Range("N44").Value = -91333
$ws.Range("H81").Value = 3278.9333
$ws.Range("J81").Value = 4359.4
$ws.Range("L81").Value = 8718.799999999999
$ws.Range("N81").Value = -10840.8
$ws.Range("H84").Value = 3278.9333
$ws.Range("J84").Value = 4359.4
$ws.Range("L84").Value = 43594
$ws.Range("N84").Value = -54202
$ws.Range("H96").Value = 5919.227
$ws.Range("I96").Value = 1701.625
$ws.Range("J96").Value = 17166.166
$ws.Range("K96").Value = 1701.625
$ws.Range("L96").Value = 17166.166
$ws.Range("M96").Value = -328.625
$ws.Range("N96").Value = -19912.166
$ws.Range("H100").Value = 633.9091
$ws.Range("I100").Value = 437.66666
$ws.Range("K100").Value = 875.33332
$ws.Range("M100").Value = -334.33332
$ws.Range("H136").Value = 3958.0679
$ws.Range("I136").Value = 3320.0408
$ws.Range("J136").Value = 7084.4
$ws.Range("K136").Value = 9960.1224
$ws.Range("L136").Value = 21253.2
$ws.Range("M136").Value = -7410.1224
$ws.Range("N136").Value = -26353.2
